$wb = $excel.ActiveWorkbook
$ws1 = $wb.Sheets("Settings")
$ws1.Range("A2").Copy()
$ws1.Range("A7").PasteSpecial(-4122, -4142, $false, $false)
